# Update slide 1 ("kathara lab" title slide) of the Kathara one-bridge lab deck:
#  - Fix the accented title text and merge it into a single run
#  - Bump the "Version" table cell from 1.1 to 2.0
#  - Add T. Caiazzi to the "Author(s)" table cell

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Title placeholder: "kathara lab" -> "Kathará lab" --------------------
$title = $s.Shapes.Item(1)
$tr = $title.TextFrame.TextRange

# The original title is split across two runs ("kathara" [misspelled, err="1"]
# and " lab"). Drop the first run's characters so only the second run (which
# already carries the desired "en-GB"/"it-IT" formatting without err="1")
# remains, then overwrite its text with the full corrected/accented title so
# the result collapses back down to a single run.
$tr.Characters(1, 7).Text = ""
$tr.Text = "Kathar" + [char]0x00E1 + " lab"

# --- Info table: Version + Author(s) --------------------------------------
$table = $s.Shapes.Item(3).Table

# Version: 1.1 -> 2.0
$table.Cell(1, 2).Shape.TextFrame.TextRange.Text = "2.0"

# Author(s): add T. Caiazzi
$table.Cell(2, 2).Shape.TextFrame.TextRange.Text = "L. Ariemma, T. Caiazzi, G. Di Battista"
